$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the "Date of Sale" values in rows 7-18 back by 30 days
# (dates move from Oct 2021 to Sep 2021 while preserving formatting).
for ($r = 7; $r -le 18; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value = $cell.Value2 - 30
}

# Update the view state: scroll position and current selection.
$ws.Range("C17:C18").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 3
